# Apply the "configured pct_start to 0.3 in params.yaml" edit to Observation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Row 20 (Change_ID 12, pct_start = 10/EPOCHS) picked up accuracy numbers
#    and a new failure message.
# ---------------------------------------------------------------------------
$ws.Range("D20").Value2 = 72.43
$ws.Range("E20").Value2 = 55.25
$ws.Range("F20").Value2 = "23rd Epoch - ValueError: Tried to step 2353 times. The specified number of total steps is 2352"

# ---------------------------------------------------------------------------
# 2. New row 21: Change_ID 13, pct_start configured to 0.3 in params.yaml.
#    Copy the look of the previous "blank result" row (row 19 / row 18 cells)
#    so the new row keeps the same yellow highlighting used for other rows.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value2 = 13
$ws.Range("B21").Value2 = "pct_start"
$ws.Range("C21").Value2 = 0.3

$ws.Range("A19").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Copy() | Out-Null
$ws.Range("C21:G21").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Rows 14-15 (Change_ID 7) switch their highlight color from yellow to
#    green, keeping the existing border/pattern structure.
# ---------------------------------------------------------------------------
$ws.Range("A14:C15").Interior.Color = 5296274
$ws.Range("D14:G15").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# 4. Update the active selection to reflect where the author finished editing.
# ---------------------------------------------------------------------------
$ws.Range("E28").Select() | Out-Null
